# Apply crypto price/volume updates (Fri May 31 22:26:19 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.640.97"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.781.93"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'595.36"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'166.47"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "3.780.56"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "'6.36"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "'36.29"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "4.418.79"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "3.778.37"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "'18.46"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").Value = "67.621.46"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'10.11"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").Value = "'457.28"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'0.0000156"
$ws.Range("E24").Value = "  +8.62%  "
$ws.Range("D25").Value = "'83.43"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'7.29"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'2.19"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "'9.22"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "3.736.32"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "'0.993"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "'45.40"
$ws.Range("E44").Value = "  +6.11%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "'47.13"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").Value = "'8.35"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").Value = "'148.31"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  -4.26%  "
$ws.Range("D50").Value = "'389.92"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("E51").Value = "  +0.93%  "
